# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet (fund-holding detail, same shape as the
#    "2021-Q4" sheet) positioned right after "2021-Q4" and before "总计".
# 2. Prepend a "2022-Q1" row to the "总计" (totals) summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet
# ---------------------------------------------------------------------
# Duplicate "2021-Q4" (same 8-column fund-holding layout / styling) and
# drop it right after the source sheet, i.e. right before "总计".
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newWs = $wb.Worksheets.Item($q4.Index + 1)
$newWs.Name = "2022-Q1"

# The template sheet carried 4 data rows, the new quarter only has 3.
$newWs.Rows.Item(5).Delete()

# Row 2 - 010783 德邦沪港深龙头混合A
$newWs.Range("A2").Value = 0
$newWs.Range("B2").NumberFormat = "@"
$newWs.Range("B2").Value = "010783"
$newWs.Range("C2").Value = "德邦沪港深龙头混合A"
$newWs.Range("D2").NumberFormat = "@"
$newWs.Range("D2").Value = "0.93"
$newWs.Range("E2").NumberFormat = "@"
$newWs.Range("E2").Value = "81.58"
$newWs.Range("F2").NumberFormat = "@"
$newWs.Range("F2").Value = "2.96"
$newWs.Range("G2").NumberFormat = "@"
$newWs.Range("G2").Value = "0.0275"
$newWs.Range("H2").Value = 9

# Row 3 - 010784 德邦沪港深龙头混合C
$newWs.Range("A3").Value = 1
$newWs.Range("B3").NumberFormat = "@"
$newWs.Range("B3").Value = "010784"
$newWs.Range("C3").Value = "德邦沪港深龙头混合C"
$newWs.Range("D3").NumberFormat = "@"
$newWs.Range("D3").Value = "0.27"
$newWs.Range("E3").NumberFormat = "@"
$newWs.Range("E3").Value = "81.58"
$newWs.Range("F3").NumberFormat = "@"
$newWs.Range("F3").Value = "2.96"
$newWs.Range("G3").NumberFormat = "@"
$newWs.Range("G3").Value = "0.0080"
$newWs.Range("H3").Value = 9

# Row 4 - 005269 华泰柏瑞港股通量化灵活配置混合
$newWs.Range("A4").Value = 2
$newWs.Range("B4").NumberFormat = "@"
$newWs.Range("B4").Value = "005269"
$newWs.Range("C4").Value = "华泰柏瑞港股通量化灵活配置混合"
$newWs.Range("D4").NumberFormat = "@"
$newWs.Range("D4").Value = "0.33"
$newWs.Range("E4").NumberFormat = "@"
$newWs.Range("E4").Value = "37.77"
$newWs.Range("F4").NumberFormat = "@"
$newWs.Range("F4").Value = "1.04"
$newWs.Range("G4").NumberFormat = "@"
$newWs.Range("G4").Value = "0.0034"
$newWs.Range("H4").Value = 5

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" sheet
# ---------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

# Shift the 3 existing data rows (2021-Q4, 2021-Q3, 2021-Q1) down by one,
# carrying formatting with them, to make room for the new row 2.
$totalWs.Range("A4:D4").Copy($totalWs.Range("A5:D5"))
$totalWs.Range("A3:D3").Copy($totalWs.Range("A4:D4"))
$totalWs.Range("A2:D2").Copy($totalWs.Range("A3:D3"))

# New first data row
$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 3
$totalWs.Range("D2").Value = 0.04

# Renumber the index column beneath it
$totalWs.Range("A3").Value = 1
$totalWs.Range("A4").Value = 2
$totalWs.Range("A5").Value = 3
